$d = $word.ActiveDocument

# 1. Update "21 years" -> "15+ years" in the professional summary
$d.Content.Find.Execute(
    "Results-driven Marketing & Data Analytics Professional with 21 years of experience",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Results-driven Marketing & Data Analytics Professional with 15+ years of experience",
    2)

# 2. Rewrite the FLEEM bullet point
$d.Content.Find.Execute(
    [char]0x2022 + " Engineered FLEEM web application using Twilio's API to make thousands of simultaneous phone calls for IVR polls",
    $true, $false, $false, $false, $false, $true, 1, $false,
    [char]0x2022 + " Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys",
    2)

# 3. Add a new bullet point after the "Developed innovative approaches..." paragraph
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Developed innovative approaches to visualizing demographic and market data, enhancing clients' understanding of research findings*") {
        $newPara = $p.Range.InsertParagraphAfter()
        $p.Next().Range.Text = [char]0x2022 + " Trained staff on building Python tooling for report generation and analysis"
        break
    }
}
